$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet to reflect the new "through" date
$ws.Name = "Through 2022-02-16"

# Update the label for the February row
$ws.Range("A3").Value = "February (through 02-16)"

# Update February row (row 3) figures
$ws.Range("D3").Value = 38
$ws.Range("E3").Value = 33
$ws.Range("F3").Value = 14
$ws.Range("G3").Value = 40
$ws.Range("H3").Value = 72
$ws.Range("I3").Value = 77

# Update Total row (row 4) figures
$ws.Range("D4").Value = 113
$ws.Range("E4").Value = 119
$ws.Range("F4").Value = 63
$ws.Range("G4").Value = 114
$ws.Range("H4").Value = 289
$ws.Range("I4").Value = 238
